$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row changes (row 1)
# ---------------------------------------------------------------------------

# D1: "Iteration 3* " -> "Iteration 3* \n(11/03/20)"
$ws.Range("D1").Value = "Iteration 3* " + [char]10 + "(11/03/20)"

# New header cells E1 ("EUC") and F1 ("SD") - same look as the other header
# cells (bold, shaded fill, wrap text) but only with a left/right border
# (no top/bottom) so they appear to be a sub-header under the main border.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Value2 = "EUC"
$ws.Range("F1").Value2 = "SD"
$ws.Range("E1:F1").Borders.Item(8).LineStyle = -4142
$ws.Range("E1:F1").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# 2. C2 / C9 become text "1(SD)" instead of the number 1
# ---------------------------------------------------------------------------
$ws.Range("C2").Value2 = "1(SD)"
$ws.Range("C9").Value2 = "1(SD)"

# ---------------------------------------------------------------------------
# 3. New column E ("Done" markers for Iteration-3 rows, styled bold +
#    centered + boxed) and column F (blank cells, "SD" marker on row 10,
#    all sharing the plain bordered/centered look used elsewhere)
# ---------------------------------------------------------------------------

# E2:E8, E20:E22 - plain existing body style (same as column C/D), left blank
$ws.Range("D2").Copy()
$ws.Range("E2:E8").PasteSpecial(-4122)
$ws.Range("E20:E22").PasteSpecial(-4122)

# E9:E19 - bold + centered + boxed style (built from the existing boxed
# style, bolded), used for the "Done" column next to Iteration-3 items
$ws.Range("C2").Copy()
$ws.Range("E9:E19").PasteSpecial(-4122)
$ws.Range("E9:E19").Font.Bold = $true

# Fill in "Done" markers on the rows that belong to Iteration 3 (D = 1)
$ws.Range("E10").Value2 = "Done"
$ws.Range("E11").Value2 = "Done"
$ws.Range("E12").Value2 = "Done"
$ws.Range("E13").Value2 = "Done"
$ws.Range("E17").Value2 = "Done"
$ws.Range("E18").Value2 = "Done"
$ws.Range("E19").Value2 = "Done"
# E9 and E14:E16 stay blank (but keep the bold boxed style applied above)

# Column F - plain body style for every row that gets an F cell
$ws.Range("D2").Copy()
$ws.Range("F2:F9").PasteSpecial(-4122)
$ws.Range("F12:F22").PasteSpecial(-4122)

# F10 - bold centered boxed style (like E10) with the "SD" marker
$ws.Range("C2").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Font.Bold = $true
$ws.Range("F10").Value2 = "SD"

# ---------------------------------------------------------------------------
# 4. View settings: zoom to 130% and move the selection to F10
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("F10").Select() | Out-Null
